$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L8").Value = "atw=at"
$ws.Range("M8").Value = 21.12

$ws.Range("L9").Value = "beta b"
$ws.Range("M9").Value = 18.37

$ws.Range("L10").Value = "bw"
$ws.Range("M10").Value = 36

$ws.Range("L11").Select()
